# The commit swaps the two fungal-observation records stored in rows 17
# and 18 (everything about the "find" itself moves to the other row),
# while the columns that are identical between the two records (C, J, K,
# N, S, T, U, V, W, Y, Z, AA, AB, AD, AE, AF, AG, AT, AW, AX, AY) are left
# untouched. Column I ("Antal") is stored as text ("10"/"6") in the
# source file, so we force text formatting before writing it back so it
# doesn't get auto-converted into a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$numericTextCols = @("I")
$plainCols = @("A", "B", "D", "E", "F", "G", "H", "P", "Q", "R")

# Force columns that hold numeric-looking text to stay text typed.
foreach ($col in $numericTextCols) {
    $ws.Range("$col" + "17:" + "$col" + "18").NumberFormat = "@"
}

# Swap the values between row 17 and row 18 for every column that differs.
foreach ($col in ($plainCols + $numericTextCols)) {
    $addr17 = "$col" + "17"
    $addr18 = "$col" + "18"
    $v17 = $ws.Range($addr17).Value2
    $v18 = $ws.Range($addr18).Value2
    $ws.Range($addr17).Value2 = $v18
    $ws.Range($addr18).Value2 = $v17
}

# Restore the default (Normal) style on the text-forced cells so no
# lingering number-format styling is left behind.
foreach ($col in $numericTextCols) {
    $ws.Range("$col" + "17:" + "$col" + "18").Style = "Normal"
}

# AC17 was empty and AC18 held "På asplåga."; after the edit it's reversed.
$ws.Range("AC17").Value2 = "På asplåga."
$ws.Range("AC18").Value2 = ""
